# "Changes of Rate Verification"
# Updates the PackageTrackNum values (col C, and col D where it mirrors C)
# for rows 2-22 to a new batch of tracking numbers, and flips the
# corresponding Rate-Verification PASS/FAIL indicator cells (cols L-W)
# to match the new run's results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tracking numbers are pure digit strings; Excel auto-detects a numeric
# literal on plain assignment and stores it as a Number. These columns
# hold identifiers (not quantities), so force them to the Text format
# first - same as typing into a cell already formatted as Text - to keep
# them stored as proper string cells (t="s"), matching the source data.
$trackCells = @(
    "C2","C3","C4","C5","D5","C6","D6","C7","D7","C8","C9","C10","C11",
    "C12","C13","D13","C14","D14","C15","D15","C16","D16","C17","D17",
    "C18","C19","C20","C21","C22"
)
foreach ($addr in $trackCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("C2").Value  = "320017965214"
$ws.Range("C3").Value  = "320017965225"
$ws.Range("C4").Value  = "320017965258"
$ws.Range("C5").Value  = "320017965270"
$ws.Range("D5").Value  = "320017965270"
$ws.Range("C6").Value  = "320017965317"
$ws.Range("D6").Value  = "320017965317"
$ws.Range("C7").Value  = "320017965339"
$ws.Range("D7").Value  = "320017965339"
$ws.Range("C8").Value  = "320017965361"
$ws.Range("C9").Value  = "320017965383"
$ws.Range("C10").Value = "320017965410"
$ws.Range("C11").Value = "320017965431"
$ws.Range("C12").Value = "320017965475"
$ws.Range("C13").Value = "320017965497"
$ws.Range("D13").Value = "320017965497"
$ws.Range("C14").Value = "320017965523"
$ws.Range("D14").Value = "320017965523"
$ws.Range("C15").Value = "320017965545"
$ws.Range("D15").Value = "320017965545"
$ws.Range("C16").Value = "320017965578"
$ws.Range("D16").Value = "320017965578"
$ws.Range("C17").Value = "320017965590"
$ws.Range("D17").Value = "320017965590"
$ws.Range("C18").Value = "320017965637"
$ws.Range("C19").Value = "320017965659"
$ws.Range("C20").Value = "320017965681"
$ws.Range("C21").Value = "320017965707"
$ws.Range("C22").Value = "320017965730"

# Rate-Verification result cells that flipped between PASS and FAIL.
$ws.Range("N4").Value  = "FAIL"
$ws.Range("P4").Value  = "FAIL"
$ws.Range("R4").Value  = "FAIL"

$ws.Range("M5").Value  = "PASS"
$ws.Range("O5").Value  = "PASS"

$ws.Range("M6").Value  = "PASS"
$ws.Range("O6").Value  = "PASS"

$ws.Range("M7").Value  = "PASS"
$ws.Range("O7").Value  = "PASS"

$ws.Range("L13").Value = "FAIL"
$ws.Range("N13").Value = "FAIL"
$ws.Range("P13").Value = "FAIL"

$ws.Range("L14").Value = "FAIL"
$ws.Range("N14").Value = "FAIL"
$ws.Range("P14").Value = "FAIL"

$ws.Range("M15").Value = "PASS"
$ws.Range("O15").Value = "PASS"
$ws.Range("W15").Value = "PASS"

$ws.Range("L16").Value = "FAIL"
$ws.Range("N16").Value = "FAIL"
$ws.Range("P16").Value = "FAIL"

$ws.Range("L17").Value = "FAIL"
$ws.Range("N17").Value = "FAIL"
$ws.Range("P17").Value = "FAIL"

$ws.Range("M19").Value = "PASS"
$ws.Range("O19").Value = "PASS"

$ws.Range("M20").Value = "PASS"
$ws.Range("O20").Value = "PASS"

$ws.Range("M21").Value = "PASS"
$ws.Range("O21").Value = "PASS"
